$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "Tiago P.-M.T"
$ws.Range("C18").Value = "[Gisele-E. D. N. D., Clesidson-Elet. Dig. Bas., Valmir-Tec. Mat. Não Metal., Aselmo-Manut. Mot. End.]"
$ws.Range("D18").Value = "[Humberto-Desenho tecnico mecanico – T1, Elcio Dec.-Desenho tecnico mecanico – T2]"
$ws.Range("E18").Value = "[Suzanny-Metalografia, -, Emerson-Comandos Eletricos, Joel L.-T. M. Metalicos]"
$ws.Range("F18").Value = "[-, Anderson-Tornearia, Victor-Ajustagem, Anderson-Metrologia 1]"

$ws.Range("B19").Value = "Tiago P.-M.T"
$ws.Range("C19").Value = "[Gisele-E. D. N. D., Clesidson-Elet. Dig. Bas., Valmir-Tec. Mat. Não Metal., Aselmo-Manut. Mot. End.]"
$ws.Range("D19").Value = "[Humberto-Desenho tecnico mecanico – T1, Elcio Dec.-Desenho tecnico mecanico – T2]"
$ws.Range("E19").Value = "[Suzanny-Metalografia, Rachel-Trat. Termicos, Emerson-Comandos Eletricos, Joel L.-T. M. Metalicos]"
$ws.Range("F19").Value = "[J. Paulo S.-Tec. Mat. Não Metal., Anderson-Tornearia, Victor-Ajustagem, Anderson-Metrologia 1]"

$ws.Range("B20").Value = "[J. Paulo S.-Tec. Mat. Não Metal., -, -, Rachel-Trat. Termicos]"
$ws.Range("C20").Value = "[Gisele-E. D. N. D., Clesidson-Elet. Dig. Bas., Valmir-Tec. Mat. Não Metal., Aselmo-Manut. Mot. End.]"
$ws.Range("D20").Value = "[Humberto-Desenho tecnico mecanico – T1, Elcio Dec.-Desenho tecnico mecanico – T2]"
$ws.Range("E20").Value = "[Suzanny-Metalografia, Rachel-Trat. Termicos, Emerson-Comandos Eletricos, Joel L.-T. M. Metalicos]"
$ws.Range("F20").Value = "[J. Paulo S.-Tec. Mat. Não Metal., Anderson-Tornearia, Victor-Ajustagem, Anderson-Metrologia 1]"

$ws.Range("B21").Value = "Gilberto-M.T.R"
$ws.Range("C21").Value = "[Gisele-E. D. N. D., Clesidson-Elet. Dig. Bas., Valmir-Tec. Mat. Não Metal., Aselmo-Manut. Mot. End.]"
$ws.Range("D21").Value = "[Humberto-Desenho tecnico mecanico – T1, Elcio Dec.-Desenho tecnico mecanico – T2]"
$ws.Range("E21").Value = "[Suzanny-Metalografia, Rachel-Trat. Termicos, Emerson-Comandos Eletricos, Joel L.-T. M. Metalicos]"
$ws.Range("F21").Value = "[J. Paulo S.-Tec. Mat. Não Metal., Anderson-Tornearia, Victor-Ajustagem, Anderson-Metrologia 1]"
